# Refresh the cryptos sheet with the latest scraped price / 1h-volume figures.
# (Some rows also gained/lost/reordered a coin entry.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings such as "211.99" or "0.0850" look numeric, so a plain
# assignment would make Excel coerce them to Double and silently drop
# formatting (e.g. the trailing zero in 0.0850). Prefixing with a literal
# apostrophe forces Excel to store the text verbatim, same as typing
# '0.0850 into the cell by hand; the apostrophe itself is not stored.
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "`'" + $text
}

$ws.Range('D2').Value = '26.226.37'
$ws.Range('E2').Value = '  +0.36%  '

$ws.Range('D3').Value = '1.588.24'
$ws.Range('E3').Value = '  +1.02%  '

Set-TextValue 'D5' '211.99'
$ws.Range('E5').Value = '  +1.80%  '

$ws.Range('E6').Value = '  +0.55%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  +0.61%  '

$ws.Range('E9').Value = '  -0.19%  '

Set-TextValue 'D10' '19.29'
$ws.Range('E10').Value = '  -1.34%  '

Set-TextValue 'D11' '0.0850'
$ws.Range('E11').Value = '  +0.70%  '

$ws.Range('D12').Value = '1.811.70'
$ws.Range('E12').Value = '  +1.01%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.588.28'
$ws.Range('E13').Value = '  +0.71%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D14' '4.01'
$ws.Range('E14').Value = '  -0.99%  '

$ws.Range('E15').Value = '  +1.30%  '

Set-TextValue 'D16' '64.26'
$ws.Range('E16').Value = '  -0.02%  '

$ws.Range('D17').Value = '26.237.76'
$ws.Range('E17').Value = '  +0.45%  '

$ws.Range('D18').Value = '0.0₃0726'
$ws.Range('E18').Value = '  +0.27%  '

Set-TextValue 'D19' '7.39'
$ws.Range('E19').Value = '  +1.55%  '

Set-TextValue 'D20' '212.81'
$ws.Range('E20').Value = '  +2.58%  '

$ws.Range('E21').Value = '  -0.04%  '

$ws.Range('E22').Value = '  +0.86%  '

Set-TextValue 'D23' '2.18'
$ws.Range('E23').Value = '  +0.22%  '

Set-TextValue 'D24' '9.01'
$ws.Range('E24').Value = '  +2.23%  '

Set-TextValue 'D25' '143.81'
$ws.Range('E25').Value = '  +0.34%  '

$ws.Range('E26').Value = '  -0.02%  '

Set-TextValue 'D27' '7.04'
$ws.Range('E27').Value = '  +1.10%  '

$ws.Range('E28').Value = '  -0.52%  '

Set-TextValue 'D29' '15.17'
$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('E30').Value = '  -1.67%  '

Set-TextValue 'D31' '1.15'
$ws.Range('E31').Value = '  +1.43%  '

$ws.Range('E32').Value = '  -0.20%  '

$ws.Range('D33').Value = '1.342.84'
$ws.Range('E33').Value = '  +5.45%  '

$ws.Range('E34').Value = '  -1.57%  '

Set-TextValue 'D35' '2.44'
$ws.Range('E35').Value = '  +0.16%  '

$ws.Range('E36').Value = '  -0.31%  '

Set-TextValue 'D37' '0.583'
$ws.Range('E37').Value = '  -4.45%  '

$ws.Range('E38').Value = '  +0.89%  '

Set-TextValue 'D39' '0.822'
$ws.Range('E39').Value = '  +1.56%  '

Set-TextValue 'D40' '5.75'
$ws.Range('E40').Value = '  +3.66%  '

$ws.Range('E41').Value = '  -0.01%  '

Set-TextValue 'D42' '0.986'
$ws.Range('E42').Value = '  -9.90%  '

Set-TextValue 'D43' '2.14'
$ws.Range('E43').Value = '  +0.61%  '

Set-TextValue 'D44' '0.767'
$ws.Range('E44').Value = '  +0.62%  '

$ws.Range('D45').Value = '1.723.79'
$ws.Range('E45').Value = '  +0.99%  '

Set-TextValue 'D46' '61.29'
$ws.Range('E46').Value = '  -1.02%  '

Set-TextValue 'D47' '85.69'
$ws.Range('E47').Value = '  -3.66%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '1.47'
$ws.Range('E48').Value = '  -2.52%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D49' '0.0975'
$ws.Range('E49').Value = '  -2.66%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.0501'
$ws.Range('E50').Value = '  -0.81%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue 'D51' '0.999'
$ws.Range('E51').Value = '  -0.14%  '

